# Add new columns I ("I0") and J ("IF") to the sheet, with header styling
# matching the existing header row (copy format from H1, the last existing
# header cell) and fill in the per-row numeric values for rows 2-75.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers -------------------------------------------------------------
# Copy H1's format (bold font, thin border, centered/top alignment) onto
# the two new header cells so they match the look of the other headers,
# then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows -------------------------------------------------------------
# Each entry is (row, I value, J value).
$data = @(
    @(2, 9, 9),
    @(3, 8, 8),
    @(4, 9, 9),
    @(5, 9, 9),
    @(6, 9, 9),
    @(7, 9, 9),
    @(8, 9, 9),
    @(9, 9, 9),
    @(10, 9, 9),
    @(11, 9, 9),
    @(12, 9, 9),
    @(13, 8, 8),
    @(14, 8, 9),
    @(15, 8, 8),
    @(16, 9, 9),
    @(17, 9, 9),
    @(18, 8, 8),
    @(19, 9, 9),
    @(20, 8, 8),
    @(21, 9, 9),
    @(22, 8, 8),
    @(23, 6, 6),
    @(24, 7, 7),
    @(25, 8, 8),
    @(26, 7, 8),
    @(27, 9, 9),
    @(28, 6, 7),
    @(29, 7, 7),
    @(30, 6, 6),
    @(31, 5, 6),
    @(32, 4, 6),
    @(33, 6, 6),
    @(34, 7, 7),
    @(35, 7, 7),
    @(36, 6, 6),
    @(37, 6, 6),
    @(38, 5, 6),
    @(39, 7, 7),
    @(40, 7, 7),
    @(41, 8, 8),
    @(42, 6, 7),
    @(43, 6, 6),
    @(44, 6, 6),
    @(45, 5, 5),
    @(46, 3, 4),
    @(47, 6, 7),
    @(48, 6, 9),
    @(49, 5, 6),
    @(50, 6, 7),
    @(51, 6, 6),
    @(52, 8, 8),
    @(53, 7, 8),
    @(54, 8, 8),
    @(55, 7, 8),
    @(56, 7, 7),
    @(57, 7, 7),
    @(58, 8, 8),
    @(59, 5, 6),
    @(60, 7, 8),
    @(61, 8, 8),
    @(62, 7, 7),
    @(63, 7, 7),
    @(64, 9, 9),
    @(65, 7, 7),
    @(66, 7, 7),
    @(67, 5, 6),
    @(68, 4, 4),
    @(69, 8, 8),
    @(70, 5, 5),
    @(71, 5, 5),
    @(72, 3, 3),
    @(73, 5, 6),
    @(74, 4, 4),
    @(75, 5, 5)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
